# Insert a new weekly price-report row for "Feria Lagunitas de Puerto Montt -
# Albahaca" at row 186. Inserting (instead of just overwriting) pushes every
# existing row from 186..211 down by one (to 187..212), which matches the
# target diff: row 186 gets brand-new data and every later row ends up
# holding what used to be in the row above it, with the former last row
# (211) surviving as the new last row (212).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 186-211 down to 187-212, leaving a blank row 186 to fill in.
$ws.Rows.Item(186).Insert()

$ws.Cells.Item(186, 1).Value  = 4
$ws.Cells.Item(186, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(186, 3).Value  = "Los Lagos"
$ws.Cells.Item(186, 4).Value  = 45218
$ws.Cells.Item(186, 5).Value  = 10
$ws.Cells.Item(186, 6).Value  = 100112052
$ws.Cells.Item(186, 7).Value  = "Albahaca"
$ws.Cells.Item(186, 8).Value  = "Sin especificar"
$ws.Cells.Item(186, 9).Value  = "Primera"
$ws.Cells.Item(186, 10).Value = 50
$ws.Cells.Item(186, 11).Value = 7000
$ws.Cells.Item(186, 12).Value = 7000
$ws.Cells.Item(186, 13).Value = 7000
$ws.Cells.Item(186, 14).Value = "`$/paquete"
$ws.Cells.Item(186, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(186, 16).Value = 7000
$ws.Cells.Item(186, 17).Value = 1
$ws.Cells.Item(186, 18).Value = "Hortaliza"
